# Adds two new columns, I (I0) and J (IF), to the data sheet, filling in
# the header labels (with the same formatting as the existing headers)
# and the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they end up with the same style (bold, bordered,
# centered) as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$i0Values = @(7,6,7,7,8,9,7,6,8,6,8,6,9,4,7,9,7,7,7,8,10,7,5,6,8,7,7,6,8,8,7,6,6,9,9,8,8,7,5,6,7,6,7,6,8,8,8,8,8,7,7,8,7,7,6,6,5,7,7,7,7,6,8,8,7,8,4,6,3,7,7,8)
$ifValues = @(7,6,8,7,8,9,7,6,8,6,8,6,9,5,7,9,7,7,7,8,10,7,5,7,8,7,8,6,8,8,7,6,7,9,9,8,8,7,5,6,8,7,7,6,8,8,8,8,8,7,7,8,7,7,6,6,5,7,7,7,7,6,8,8,8,8,4,6,4,7,7,8)

for ($i = 0; $i -lt $i0Values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$i]
    $ws.Cells.Item($row, 10).Value = $ifValues[$i]
}
